$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Bump the "updated" date (column C, rows 2-221) from 45189 to 45190 (one day later),
# matching the stored numeric date serial values exactly.
for ($r = 2; $r -le 221; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
